$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.748.14"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +4.68%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'2.248.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +3.77%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = "'  +0.23%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'248.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.81%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("E6").Value = "'  +2.73%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'70.27"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +5.74%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("E8").Value = "'  +0.06%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.662"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +16.39%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'38.93"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +9.47%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'59.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +1.56%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'0.0962"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +4.24%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'7.45"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +8.07%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("E14").Value = "'  +0.50%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'2.576.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.63%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'14.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +3.72%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'0.876"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.83%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'2.277.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +5.43%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'42.689.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +4.66%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'0.0₃0982"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +4.80%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'6.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +3.10%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'72.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +2.10%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'234.94"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.73%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = "'  -3.13%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("E25").Value = "'  +6.12%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("E26").Value = "'  +0.05%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'11.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -1.74%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'2.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.05%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'3.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -1.47%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'2.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +9.99%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'167.53"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.74%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'20.80"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +3.21%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = "'  +13.98%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = "'  +5.08%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.0798"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +7.74%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'31.23"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +26.34%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  +4.42%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'4.43"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +12.14%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'4.68"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +3.25%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  +8.44%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'2.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +5.52%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'12.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +7.17%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D44").Value = "'62.11"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.01%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").Value = "'FTXToken"
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = "'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").Value = "'4.90"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.94%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").Value = "'Algorand"
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").Value = "'0.201"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +4.83%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").Value = "'FraxShare"
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").Value = "'8.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.86%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  +2.50%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.998"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.23%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("E50").Value = "'  +2.48%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  +3.55%  "
$ws.Range("E51").Style = "Normal"
